# Update "countries & provincias Spain" data dump.
# - refreshed timestamp string
# - three pairs of countries swap rank (same row, new country name + new stats)
# - refreshed case/death counters for a batch of existing country rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 28 de Septiembre de 2020 a las 19:11"

# --- Country name swaps (rank changed, row stays, name + numbers update) ---
$ws.Range("A62").Value = "Argelia"
$ws.Range("A63").Value = "Moldavia"

$ws.Range("A97").Value = "Namibia"
$ws.Range("A98").Value = "Malasia"

$ws.Range("A207").Value = "Santa Lucia"
$ws.Range("A208").Value = "Timor Oriental"

# --- Numeric refreshes: Casos totales / Nuevos casos / Casos activos /
#     Recuperados / Casos criticos / Muertes hoy / Muertes ---

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 7328926
$ws.Range("C4").Value = 7583
$ws.Range("D4").Value = 4577350
$ws.Range("E4").Value = 2542050
$ws.Range("G4").Value = 73
$ws.Range("H4").Value = 209526

# Row 6 - Brasil
$ws.Range("B6").Value = 4736831
$ws.Range("C6").Value = 4522
$ws.Range("E6").Value = 534967

# Row 15 - Chile
$ws.Range("B15").Value = 459671
$ws.Range("C15").Value = 1770
$ws.Range("D15").Value = 433016
$ws.Range("E15").Value = 13957
$ws.Range("G15").Value = 57
$ws.Range("H15").Value = 12698

# Row 21 - Turquia
$ws.Range("B21").Value = 315845
$ws.Range("C21").Value = 1412
$ws.Range("D21").Value = 277052
$ws.Range("E21").Value = 30731
$ws.Range("G21").Value = 65
$ws.Range("H21").Value = 8062

# Row 25 - Alemania
$ws.Range("B25").Value = 287786
$ws.Range("C25").Value = 1448
$ws.Range("E25").Value = 27447

# Row 27 - Israel
$ws.Range("B27").Value = 233118
$ws.Range("C27").Value = 2092
$ws.Range("D27").Value = 164980
$ws.Range("E27").Value = 66639
$ws.Range("G27").Value = 33
$ws.Range("H27").Value = 1499

# Row 57 - Chequia
$ws.Range("B57").Value = 65313
$ws.Range("C57").Value = 716
$ws.Range("E57").Value = 33430
$ws.Range("G57").Value = 9
$ws.Range("H57").Value = 615

# Row 62 - Argelia (now; this row's stats)
$ws.Range("B62").Value = 51213
$ws.Range("C62").Value = 146
$ws.Range("D62").Value = 35962
$ws.Range("E62").Value = 13532
$ws.Range("G62").Value = 5
$ws.Range("H62").Value = 1719

# Row 63 - Moldavia (now; this row's stats)
$ws.Range("B63").Value = 51194
$ws.Range("C63").Value = 319
$ws.Range("D63").Value = 38217
$ws.Range("E63").Value = 11676
$ws.Range("G63").Value = 14
$ws.Range("H63").Value = 1301

# Row 72 - Kenia
$ws.Range("B72").Value = 38168
$ws.Range("C72").Value = 53
$ws.Range("D72").Value = 24681
$ws.Range("E72").Value = 12787
$ws.Range("G72").Value = 9
$ws.Range("H72").Value = 700

# Row 88 - Madagascar
$ws.Range("B88").Value = 16348
$ws.Range("C88").Value = 63
$ws.Range("D88").Value = 14947
$ws.Range("E88").Value = 1172

# Row 97 - Namibia (now; this row's stats)
$ws.Range("B97").Value = 11121
$ws.Range("C97").Value = 88
$ws.Range("D97").Value = 8787
$ws.Range("E97").Value = 2213
$ws.Range("G97").Value = 1
$ws.Range("H97").Value = 121

# Row 98 - Malasia (now; this row's stats)
$ws.Range("B98").Value = 11034
$ws.Range("C98").Value = 115
$ws.Range("D98").Value = 9889
$ws.Range("E98").Value = 1011
$ws.Range("H98").Value = 134

# Row 111 - Mozambique
$ws.Range("B111").Value = 8288
$ws.Range("C111").Value = 305
$ws.Range("D111").Value = 4836
$ws.Range("E111").Value = 3393
$ws.Range("G111").Value = 1
$ws.Range("H111").Value = 59

# Row 121 - Republica de Yibuti
$ws.Range("B121").Value = 5410
$ws.Range("C121").Value = 1
$ws.Range("E121").Value = 9

# Row 142 - Sri Lanka
$ws.Range("B142").Value = 3362
$ws.Range("C142").Value = 2
$ws.Range("E142").Value = 139

# Row 148 - Sudan del Sur
$ws.Range("B148").Value = 2692
$ws.Range("C148").Value = 6
$ws.Range("E148").Value = 1353

# Row 153 - Yemen
$ws.Range("B153").Value = 2031
$ws.Range("C153").Value = 1
$ws.Range("E153").Value = 178

# Row 184 - Curazao
$ws.Range("B184").Value = 364
$ws.Range("C184").Value = 4
$ws.Range("E184").Value = 215

# Row 197 - Bonaire, San Eustaquio y Saba
$ws.Range("B197").Value = 88
$ws.Range("C197").Value = 3
$ws.Range("E197").Value = 66
